# Update the "想去人数" (F column) values for the affected rows on both the
# "展览" and "全部类型" worksheets. These two sheets carry duplicated rows of
# data, and the same set of rows/values were refreshed in the source data.

$wb = $excel.ActiveWorkbook

# Map of row number -> new F-column value (applies to both sheets below).
$updates = @{
    5  = 4698
    7  = 412
    8  = 1420
    9  = 928
    10 = 57
    11 = 1235
    13 = 765
    15 = 64
    16 = 29
    18 = 25
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
